$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53; everything below (old rows 53-96) shifts
# down to 54-97, preserving all of their existing values/styles.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new weekly record. The
# non-numeric / unchanged fields mirror the row that used to occupy this
# slot (same market, region, category, etc.) per the diff.
$ws.Range("A53").Value = 11
$ws.Range("B53").Value = "Vega Monumental Concepción"
$ws.Range("C53").Value = "Bíobío"
$ws.Range("D53").Value = 44741
$ws.Range("E53").Value = 8
$ws.Range("F53").Value = 100112001
$ws.Range("G53").Value = "Berenjena"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 100
$ws.Range("K53").Value = 12000
$ws.Range("L53").Value = 14000
$ws.Range("M53").Value = 13000
$ws.Range("N53").Value = "$/caja 60 unidades"
$ws.Range("O53").Value = "Región de Arica y Parinacota"
$ws.Range("P53").Value = 217
$ws.Range("Q53").Value = 60
$ws.Range("R53").Value = "Hortaliza"
